$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "02/20/2026"
$ws.Range("B88").Value = 9017.030000000001
$ws.Range("C88").Value = 0.2497591506634206
$ws.Range("D88").Value = 0.7502408493365794
$ws.Range("E88").Value = -353.1
$ws.Range("F88").Value = -37.76
$ws.Range("G88").Value = -24226.81
$ws.Range("H88").Value = -78.17
$ws.Range("I88").Value = -1150.76
$ws.Range("J88").Value = -33.82
$ws.Range("K88").Value = -25377.57
$ws.Range("L88").Value = -73.78
